$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove all existing hyperlinks on the sheet (Hyperlinks.Delete clears the whole sheet in this engine)
$ws.Range("A1").Hyperlinks.Delete()

# Row 2
$ws.Range("A2").Value = "2025-11-29 12:34:15"
$ws.Range("B2").Value = "【急募】フロントエンド開発者募集!React/TypeScriptでのシステム構築"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5443491"
$ws.Range("G2").Value = 323
$ws.Range("H2").Value = "🔥React,TypeScript ◆開発"

# Row 3
$ws.Range("A3").Value = "2025-11-29 12:34:15"
$ws.Range("B3").Value = "【急募】在庫・販売管理ツールの開発依頼"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5443889"
$ws.Range("G3").Value = 170
$ws.Range("H3").Value = "◆ツール,開発 ◇管理"

# Row 4
$ws.Range("A4").Value = "2025-11-29 12:34:15"
$ws.Range("B4").Value = "マンション管理組合のシステム設計構築依頼"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5443592"
$ws.Range("G4").Value = 60
$ws.Range("H4").Value = "◇管理"

# Row 5
$ws.Range("A5").Value = "2025-11-29 12:34:15"
$ws.Range("B5").Value = "【Apache Answer構築】弁護士ドットコムのような専門家Q&Aサイトのサーバー構築・初期設定"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5443617"
$ws.Range("G5").Value = 38
$ws.Range("H5").Value = "◇サイト"

# Row 6
$ws.Range("A6").Value = "2025-11-29 12:34:15"
$ws.Range("B6").Value = "【急募】Wartalesの武器アイコンとモデルを日本刀に差し替え"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5443568"
$ws.Range("G6").Value = 13

# Re-create hyperlinks for F2:F6 in order, then restore the Hyperlink cell style
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5443491")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5443889")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5443592")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5443617")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5443568")
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 6).Style = "Hyperlink"
}
